$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: replace Forename/Surname placeholders with real values
$ws.Range("A2").Value = "PlanIt"
$ws.Range("B2").Value = "Solution"

# Row 3: new row of data
$ws.Range("A3").Value = "PlanIt"
$ws.Range("B3").Value = "Solution"
$ws.Range("C3").Value = 123
$ws.Range("D3").Value = "Asd"
$ws.Range("E3").Value = "message"

# Update the active selection to match the authored workbook state
$ws.Range("E4").Select()
